# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by copying "2021-Q3" (so it keeps
#    the same look & feel / cell styles), placed right after it.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2021-Q3")
$q3.Copy($null, $q3)
$q1_2022 = $wb.Worksheets.Item("2021-Q3 (2)")
$q1_2022.Name = "2022-Q1"

# Header row - fix the column D label ("基金金额" -> "基金规模"),
# the rest of the header row is identical to "2021-Q3".
$q1_2022.Range("B1").Value = "基金代码"
$q1_2022.Range("C1").Value = "基金名称"
$q1_2022.Range("D1").Value = "基金规模"
$q1_2022.Range("E1").Value = "股票总仓位"
$q1_2022.Range("F1").Value = "仓位占比"
$q1_2022.Range("G1").Value = "持有市值(亿元)"
$q1_2022.Range("H1").Value = "仓位排名"

# Data row - the leading apostrophe forces Excel to store the
# numeric-looking values (fund code, percentages, ...) as text,
# matching the source data's text columns.
$q1_2022.Range("B2").Value = "'233009"
$q1_2022.Range("C2").Value = "大摩多因子精选策略混合"
$q1_2022.Range("D2").Value = "'6.77"
$q1_2022.Range("E2").Value = "'89.73"
$q1_2022.Range("F2").Value = "'1.18"
$q1_2022.Range("G2").Value = "'0.0799"
$q1_2022.Range("H2").Value = 5

# ------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: a new 2022-Q1 row is inserted
#    on top, the existing rows shift down and get renumbered.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Grab the format of an existing numbered row so the freshly
# created row 4 can match the existing style (s="2" on column A).
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q1"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.1

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.09

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.08

$wb.Worksheets.Item("2021-Q1").Select()
